$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97, shifting the existing row 97 (and everything
# below it) down by one. This matches the diff: all rows from the old 97..134
# move to 98..135, and the sheet's used range grows from A1:R134 to A1:R135.
$ws.Range("A97").EntireRow.Insert()

# Populate the freshly inserted row 97 with the new record's data. The
# non-numeric/categorical columns (market, region, product, quality, unit,
# origin, kg/units, classification) repeat the same constant values used by
# every other row in this block.
$ws.Cells.Item(97, 1).Value = 11
$ws.Cells.Item(97, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(97, 3).Value = "Bíobío"
$ws.Cells.Item(97, 4).Value = 44726
$ws.Cells.Item(97, 5).Value = 8
$ws.Cells.Item(97, 6).Value = 100112043
$ws.Cells.Item(97, 7).Value = "Pepino ensalada"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 270
$ws.Cells.Item(97, 11).Value = 17000
$ws.Cells.Item(97, 12).Value = 18000
$ws.Cells.Item(97, 13).Value = 17444
$ws.Cells.Item(97, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(97, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(97, 16).Value = 291
$ws.Cells.Item(97, 17).Value = 60
$ws.Cells.Item(97, 18).Value = "Hortaliza"
